# Adjust membership function of blink rate.
# Blink rate now ranges 0..21 (was 1..20) and, per the new fuzzy grading
# table, blink rate and its membership value are equivalent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): blink-rate axis labels, now 0..21 across B1:W1 (was 1..20 across B1:U1)
$headerValues = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
for ($c = 0; $c -lt $headerValues.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $headerValues[$c]
}

# Membership-function grid (rows 2-12, columns B-W): blink rate is now equivalent to its membership value
$gradeRows = @(
    @(0.59,0.59,0.59,0.59,0.59,0.59,0.59,0.59,0.59,0.55,0.52,0.47,0.42,0.37,0.33,0.28,0.28,0.26,0.24,0.2,0.12,0),
    @(0.6,0.6,0.6,0.6,0.6,0.6,0.6,0.6,0.6,0.57,0.53,0.48,0.43,0.38,0.34,0.29,0.29,0.27,0.25,0.21,0.13,0),
    @(0.64,0.64,0.64,0.64,0.64,0.64,0.64,0.64,0.64,0.61,0.57,0.52,0.46,0.41,0.37,0.32,0.31,0.3,0.28,0.23,0.13,0),
    @(0.71,0.71,0.71,0.71,0.71,0.71,0.71,0.71,0.71,0.67,0.63,0.58,0.52,0.46,0.41,0.36,0.36,0.35,0.31,0.23,0.13,0),
    @(0.82,0.82,0.82,0.82,0.82,0.82,0.82,0.82,0.82,0.78,0.71,0.63,0.5600000000000001,0.5,0.45,0.4,0.39,0.36,0.31,0.23,0.13,0),
    @(1,1,1,1,1,1,1,1,1,0.84,0.73,0.65,0.58,0.52,0.47,0.41,0.4,0.36,0.31,0.23,0.13,0),
    @(1,1,1,1,1,1,1,1,1,0.84,0.73,0.65,0.58,0.52,0.49,0.49,0.48,0.44,0.4,0.34,0.27,0.25),
    @(1,1,1,1,1,1,1,1,1,0.84,0.73,0.65,0.58,0.5600000000000001,0.55,0.55,0.54,0.52,0.48,0.44,0.41,0.39),
    @(1,1,1,1,1,1,1,1,1,0.84,0.73,0.65,0.63,0.62,0.61,0.61,0.6,0.58,0.55,0.51,0.49,0.48),
    @(1,1,1,1,1,1,1,1,1,0.84,0.75,0.7,0.67,0.66,0.65,0.65,0.64,0.62,0.6,0.57,0.55,0.54),
    @(1,1,1,1,1,1,1,1,1,0.86,0.78,0.74,0.71,0.6899999999999999,0.68,0.68,0.67,0.66,0.64,0.61,0.59,0.59)
)
for ($r = 0; $r -lt $gradeRows.Length; $r++) {
    $rowValues = $gradeRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowValues[$c]
    }
}
